$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Account-number column (C) holds 24-digit strings that must stay text, not be
# coerced to floating point numbers.
$ws.Range("C2:C5").NumberFormat = "@"

# Row 2 -> NABIL KAMAL
$ws.Range("A2").Value = "NABIL KAMAL"
$ws.Range("B2").Value = "L3578354"
$ws.Range("C2").Value = "345534544587485743558673"
$ws.Range("D2").Value = "AGG1"
$ws.Range("E2").Value = "BP"
$ws.Range("F2").Value = "Point de vente"
$ws.Range("G2").Value = "990/PV 01"
$ws.Range("H2").Value = "trimestrielle"
$ws.Range("I2").Value = 10000
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 9000

# Row 3 -> KHADIJA LALA
$ws.Range("A3").Value = "KHADIJA LALA"
$ws.Range("B3").Value = "K5443645"
$ws.Range("C3").Value = "354564564324158786713544"
$ws.Range("D3").Value = "AG 100"
$ws.Range("E3").Value = "BP"
$ws.Range("F3").Value = "Logement de fonction"
$ws.Range("G3").Value = "044/LF/FES VILLE /AV1"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 10000
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 8500

# Row 4 -> SAMIRA TATA
$ws.Range("A4").Value = "SAMIRA TATA"
$ws.Range("B4").Value = "D524564"
$ws.Range("C4").Value = "335463513748543615567464"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "BP"
$ws.Range("F4").Value = "Supervision"
$ws.Range("G4").Value = "554/SUP FES 1"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 10000
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 8500

# Row 5 -> KARIM JALAL
$ws.Range("A5").Value = "KARIM JALAL"
$ws.Range("B5").Value = "P5874857"
$ws.Range("C5").Value = "548748641684867461687153"
$ws.Range("D5").Value = "FES SUD"
$ws.Range("E5").Value = "BMCE"
$ws.Range("F5").Value = "Point de vente"
$ws.Range("G5").Value = "800/PV FES 1"
$ws.Range("H5").Value = "mensuelle"
$ws.Range("I5").Value = 10000
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 8500

# Row 6 -> blank totals row
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("I6").Value = 40000
$ws.Range("J6").Value = 5500
$ws.Range("K6").Value = 34500

# Remove old rows 7 and 8 (LATIFA FIFA entry + old totals row)
$ws.Range("A7:K8").EntireRow.Delete()
